$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the "Heat map on international sales..." / "1. Web scraping on 25
#    movies?" row (old row 17). Everything below shifts up by one row.
# ---------------------------------------------------------------------------
$ws.Rows(17).Delete()

# ---------------------------------------------------------------------------
# 2. Status ("E" column) updates for the ML work-stream (rows shifted the
#    same as before the delete, since these rows are all above row 17).
# ---------------------------------------------------------------------------
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 1

# ---------------------------------------------------------------------------
# 3. "Vote on James Bond next actor" row (row 15): drop the trailing "?" in
#    the output title, extend the work-involved note with a second line, and
#    bump the status / row height.
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = "Vote on James Bond next actor"
$ws.Range("B15").Value = "1. New python/JS library that does survey`n2. Create a code that connects to database and update the database accordingly (stretch goal)"
$ws.Range("E15").Value = 0.9
$ws.Rows(15).RowHeight = 45

# ---------------------------------------------------------------------------
# 4. "Using Flask & D3 visualisation on Bond Girls" row -- now row 18 after
#    the deletion above -- gets a status value of 1 (was blank).
# ---------------------------------------------------------------------------
$ws.Range("E18").Value = 1

# ---------------------------------------------------------------------------
# 5. "Other" section status updates -- rows 23-25 after the deletion above
#    (formerly rows 24-26).
# ---------------------------------------------------------------------------
$ws.Range("E23").Value = 0.8
$ws.Range("E24").Value = 1
$ws.Range("E25").Value = 0.9

# ---------------------------------------------------------------------------
# 6. Refresh the AutoFilter so its range matches the new used range (A1:F26)
#    instead of the stale A1:F27, and update the workbook-level
#    _FilterDatabase defined name to match.
# ---------------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:F26").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name() -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$F`$26"
    }
}

# ---------------------------------------------------------------------------
# 7. Move the active selection to C11:C14 (active cell C11), matching the
#    author's cursor position when they saved.
# ---------------------------------------------------------------------------
$ws.Range("C11:C14").Select()
